$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")
$tbl = $ws.ListObjects.Item("Data_matricen")

# Insert 3 new rows right below the first data row (row 2) and fill them in,
# shifting all existing observations down by three rows.
$ws.Rows("3:5").Insert()

$ws.Range("A3").Value = "Anker"
$ws.Range("B3").Value = "Jon"
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = "*"

$ws.Range("A4").Value = "Bo"
$ws.Range("B4").Value = "Immanuel"
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = "*"

$ws.Range("A5").Value = "Hugo"
$ws.Range("B5").Value = "Børge"
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = "*"

# Grow the table to cover the new rows (A1:H31).
$tbl.Resize($ws.Range("A1:H31"))

# Fill in the calculated columns for the new rows, and re-apply them across
# the rest of the table too: resizing the table can otherwise leave stale /
# mis-anchored formulas in the last few rows of the (previous) table range.
$ws.Range("E3:E31").Formula = "=UPPER(RIGHT(Data_matricen[[#This Row],[Navn1]],1))"
$ws.Range("F3:F31").Formula = "=LEFT(Data_matricen[[#This Row],[Navn2]],1)"

# Re-enter Navn1_godt?/Navn2_godt? for every data row as real boolean values
# (previously these were the text strings "true"/"false").
$ws.Range("G2").Value = $false
$ws.Range("H2").Value = $false

$ws.Range("G3").Value = $true
$ws.Range("H3").Value = $true

$ws.Range("G4").Value = $true
$ws.Range("H4").Value = $true

$ws.Range("G5").Value = $true
$ws.Range("H5").Value = $true

$ws.Range("G6").Value = $false
$ws.Range("H6").Value = $false

$ws.Range("G7").Value = $false
$ws.Range("H7").Value = $true

$ws.Range("G8").Value = $false
$ws.Range("H8").Value = $false

$ws.Range("G9").Value = $false
$ws.Range("H9").Value = $false

$ws.Range("G10").Value = $false
$ws.Range("H10").Value = $true

$ws.Range("G11").Value = $false
$ws.Range("H11").Value = $false

$ws.Range("G12").Value = $false
$ws.Range("H12").Value = $false

$ws.Range("G13").Value = $false
$ws.Range("H13").Value = $false

$ws.Range("G14").Value = $false
$ws.Range("H14").Value = $false

$ws.Range("G15").Value = $false
$ws.Range("H15").Value = $true

$ws.Range("G16").Value = $false
$ws.Range("H16").Value = $false

$ws.Range("G17").Value = $false
$ws.Range("H17").Value = $true

$ws.Range("G18").Value = $true
$ws.Range("H18").Value = $true

$ws.Range("G19").Value = $false
$ws.Range("H19").Value = $false

$ws.Range("G20").Value = $true
$ws.Range("H20").Value = $false

$ws.Range("G21").Value = $false
$ws.Range("H21").Value = $true

$ws.Range("G22").Value = $false
$ws.Range("H22").Value = $false

$ws.Range("G23").Value = $false
$ws.Range("H23").Value = $false

$ws.Range("G24").Value = $true
$ws.Range("H24").Value = $true

$ws.Range("G25").Value = $false
$ws.Range("H25").Value = $false

$ws.Range("G26").Value = $false
$ws.Range("H26").Value = $false

$ws.Range("G27").Value = $false
$ws.Range("H27").Value = $false

$ws.Range("G28").Value = $false
$ws.Range("H28").Value = $false

$ws.Range("G29").Value = $false
$ws.Range("H29").Value = $false

$ws.Range("G30").Value = $true
$ws.Range("H30").Value = $true

$ws.Range("G31").Value = $true
$ws.Range("H31").Value = $true

# Reflect the author's final cursor position on the sheet.
$ws.Range("E40").Select() | Out-Null

$wb.Save()
